$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged values, set for safety)
$ws.Range("A1").Value = "individual"
$ws.Range("B1").Value = "Birth_Month"
$ws.Range("C1").Value = "Height_inches"

# Data rows 2-33: individual id, Birth_Month, Height_inches
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "January"
$ws.Range("C2").Value = 70
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "September"
$ws.Range("C3").Value = 64
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "March"
$ws.Range("C4").Value = 72
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "April"
$ws.Range("C5").Value = 61
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "April"
$ws.Range("C6").Value = 55
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "October"
$ws.Range("C7").Value = 65
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "November"
$ws.Range("C8").Value = 72
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "January"
$ws.Range("C9").Value = 75
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "March"
$ws.Range("C10").Value = 69
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "June"
$ws.Range("C11").Value = 75
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "September"
$ws.Range("C12").Value = 76
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "December"
$ws.Range("C13").Value = 70
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "December"
$ws.Range("C14").Value = 70
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "April"
$ws.Range("C15").Value = 69
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "February"
$ws.Range("C16").Value = 69
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "July"
$ws.Range("C17").Value = 65
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "August"
$ws.Range("C18").Value = 65
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "July"
$ws.Range("C19").Value = 64
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "August"
$ws.Range("C20").Value = 58
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "May"
$ws.Range("C21").Value = 57
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "June"
$ws.Range("C22").Value = 64
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "October"
$ws.Range("C23").Value = 59
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "November"
$ws.Range("C24").Value = 59
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "December"
$ws.Range("C25").Value = 63
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "January"
$ws.Range("C26").Value = 64
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "May"
$ws.Range("C27").Value = 66
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "May"
$ws.Range("C28").Value = 65
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "November"
$ws.Range("C29").Value = 67
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "April"
$ws.Range("C30").Value = 69
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "May"
$ws.Range("C31").Value = 72
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "June"
$ws.Range("C32").Value = 70
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "November"
$ws.Range("C33").Value = 71

# Remove old trailing rows (34-53), shrinking the used range to A1:C33
$ws.Range("A34:C53").ClearContents()

# Update the view: scroll so row 25 is at top, and select A31 (per saved view state)
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("A31").Select()
